# Updated cryptos list on Sat Jul 27 09:52:50 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: some Price values (column D) look like plain decimal numbers
# (e.g. "588.48"). Excel's Range.Value setter auto-coerces such strings to
# numeric cells, which does not match the source data (plain text cells).
# Forcing the NumberFormat to Text ("@") before assignment keeps the value
# as a string; resetting the Style to "Normal" afterwards removes the
# extra style index that the NumberFormat change would otherwise leave
# behind, so the cell ends up with no explicit style - matching the
# original workbook's cells.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "68.172.79"
$ws.Range("E2").Value = "  +1.17%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.275.72"
$ws.Range("E3").Value = "  +0.47%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.02%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "588.48"
$ws.Range("E5").Value = "  +1.71%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "186.36"
$ws.Range("E6").Value = "  +3.63%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.00%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  +0.01%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  +4.21%  "

# Row 10 - Toncoin
$ws.Range("E10").Value = "  -0.04%  "

# Row 11 - Cardano
$ws.Range("E11").Value = "  +0.53%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "3.842.71"
$ws.Range("E12").Value = "  +0.49%  "

# Row 14 - Avalanche
Set-TextValue $ws.Range("D14") "28.68"
$ws.Range("E14").Value = "  +1.62%  "

# Row 15 - WrappedBTC
$ws.Range("D15").Value = "68.166.20"
$ws.Range("E15").Value = "  +1.23%  "

# Row 16 - ShibaInu
$ws.Range("E16").Value = "  +2.22%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "3.280.58"
$ws.Range("E17").Value = "  +0.76%  "

# Row 18 - Polkadot
$ws.Range("E18").Value = "  -0.48%  "

# Row 19 - Chainlink
$ws.Range("E19").Value = "  +1.67%  "

# Row 20 - BitcoinCash
Set-TextValue $ws.Range("D20") "382.05"
$ws.Range("E20").Value = "  +1.14%  "

# Row 21 - Uniswap
Set-TextValue $ws.Range("D21") "7.74"
$ws.Range("E21").Value = "  +1.15%  "

# Row 22 - Dai
$ws.Range("E22").Value = "  +0.02%  "

# Row 23 - Litecoin
Set-TextValue $ws.Range("D23") "71.54"
$ws.Range("E23").Value = "  +0.22%  "

# Row 24 - was Polygon, now PEPE
$ws.Range("B24").Value = "PEPE"
$ws.Range("C24").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue $ws.Range("D24") "0.0000121"
$ws.Range("E24").Value = "  +1.96%  "

# Row 25 - was PEPE, now Polygon
$ws.Range("B25").Value = "Polygon"
$ws.Range("C25").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws.Range("D25") "0.515"
$ws.Range("E25").Value = "  +0.32%  "

# Row 26 - Kaspa
$ws.Range("E26").Value = "  +6.08%  "

# Row 27 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D27") "9.78"
$ws.Range("E27").Value = "  -1.26%  "

# Row 28 - Binance-PegBSC-USD
$ws.Range("E28").Value = "  +0.02%  "

# Row 29 - NEARProtocol
$ws.Range("E29").Value = "  +2.76%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +0.85%  "

# Row 31 - EthereumClassic
Set-TextValue $ws.Range("D31") "22.96"
$ws.Range("E31").Value = "  +1.54%  "

# Row 32 - Aptos
$ws.Range("E32").Value = "  +4.87%  "

# Row 33 - Fetch.AI
$ws.Range("E33").Value = "  +0.69%  "

# Row 34 - USDe
$ws.Range("E34").Value = "  +0.02%  "

# Row 35 - ImmutableX
$ws.Range("E35").Value = "  +2.27%  "

# Row 36 - Monero
Set-TextValue $ws.Range("D36") "163.12"
$ws.Range("E36").Value = "  -0.52%  "

# Row 37 - Stacks
$ws.Range("E37").Value = "  +0.05%  "

# Row 38 - Mantle
$ws.Range("E38").Value = "  -2.02%  "

# Row 39 - RenderToken
Set-TextValue $ws.Range("D39") "6.81"
$ws.Range("E39").Value = "  +1.74%  "

# Row 40 - EnergySwap
Set-TextValue $ws.Range("D40") "26.63"
$ws.Range("E40").Value = "  -1.53%  "

# Row 41 - was Filecoin, now dogwifhat
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws.Range("D41") "2.66"
$ws.Range("E41").Value = "  +0.65%  "

# Row 42 - was dogwifhat, now Filecoin
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D42") "4.62"
$ws.Range("E42").Value = "  +4.82%  "

# Row 43 - OKB
Set-TextValue $ws.Range("D43") "41.35"
$ws.Range("E43").Value = "  +2.13%  "

# Row 44 - Hedera
$ws.Range("E44").Value = "  +2.34%  "

# Row 45 - InjectiveProtocol
Set-TextValue $ws.Range("D45") "25.52"
$ws.Range("E45").Value = "  -1.30%  "

# Row 46 - was Bittensor, now Maker
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "2.645.01"
$ws.Range("E46").Value = "  -4.70%  "

# Row 47 - was Maker, now Bittensor
$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D47") "344.80"
$ws.Range("E47").Value = "  -0.54%  "

# Row 48 - VeChain
$ws.Range("E48").Value = "  +1.50%  "

# Row 49 - Arweave
Set-TextValue $ws.Range("D49") "32.02"
$ws.Range("E49").Value = "  +3.54%  "

# Row 50 - ONDO
$ws.Range("E50").Value = "  +1.08%  "

# Row 51 - Stellar
$ws.Range("E51").Value = "  -0.14%  "
